# Fix Training Data Issue (#48)
# The "Date" column (BF) held values formatted like "4-19-2013-14" (a
# leftover concatenation of the game date and season). These need to be
# corrected to a proper ISO-style date string "2014-04-19" (the data was
# off by one day due to how NBA stats were displayed), while remaining a
# plain text value (Excel would otherwise happily reinterpret a string
# like "2014-04-19" as a real date serial number, which is not what the
# source data represents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "2014-04-19"

# Column BF, rows 2-31 all contain the same stale date string
# ("4-19-2013-14") and need to become "2014-04-19".
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    # Assign via a text formula and then paste the computed result back
    # as a literal value. This keeps the cell a plain text string
    # ("2014-04-19") instead of Excel auto-converting the literal into
    # a date serial number, and it does so without layering a new
    # number-format style onto the cell (matching the original
    # unformatted cell).
    $cell.Formula = '="' + $newDate + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
